$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "Solar" column (E) values for 2022 (row 24) and 2024 (row 26)
# to reflect updated upstream data through 2024.
$ws.Range("E24").Value = 2
$ws.Range("E26").Value = 8
